$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.05840758140714978
$ws.Range("C2").Value = 0.4409791216474546
$ws.Range("D2").Value = 0.2667866710384161
$ws.Range("E2").Value = 0.516513960158306
$ws.Range("F2").Value = 0.5247350061034883
$ws.Range("B3").Value = 0.04755880890717359
$ws.Range("C3").Value = 0.5329826934417879
$ws.Range("D3").Value = 0.5613959449393942
$ws.Range("E3").Value = 0.749263601771362
$ws.Range("F3").Value = 0.7645582327128724
$ws.Range("B4").Value = 0.5041968073251691
$ws.Range("C4").Value = 0.8176875864659832
$ws.Range("D4").Value = 4.040021400480529
$ws.Range("E4").Value = 2.009980447785632
$ws.Range("F4").Value = 1.98944442380694
$ws.Range("B5").Value = 0.1550513310109709
$ws.Range("C5").Value = 1.247131068598841
$ws.Range("D5").Value = 7.211909333362784
$ws.Range("E5").Value = 2.685499829335832
$ws.Range("F5").Value = 2.741275197758698
$ws.Range("B6").Value = 0.06215009439971476
$ws.Range("C6").Value = 1.169868415477642
$ws.Range("D6").Value = 7.120640959472445
$ws.Range("E6").Value = 2.668452914981347
$ws.Range("F6").Value = 2.727685511402118
$ws.Range("B7").Value = 0.1800417658608053
$ws.Range("C7").Value = 1.300679498311106
$ws.Range("D7").Value = 7.248419567119259
$ws.Range("E7").Value = 2.692288908553326
$ws.Range("F7").Value = 2.7466351660821
$ws.Range("B8").Value = 0.01362888296484679
$ws.Range("C8").Value = 1.27992735595733
$ws.Range("D8").Value = 7.348708252423751
$ws.Range("E8").Value = 2.710850097741251
$ws.Range("F8").Value = 2.771740654578703
$ws.Range("B9").Value = 0.08542880453419803
$ws.Range("C9").Value = 1.381444734124759
$ws.Range("D9").Value = 7.473256470929509
$ws.Range("E9").Value = 2.733725749033635
$ws.Range("F9").Value = 2.793800303403652
$ws.Range("B10").Value = 0.007591499017025043
$ws.Range("C10").Value = 1.301026266972584
$ws.Range("D10").Value = 7.361617303534813
$ws.Range("E10").Value = 2.713230049873179
$ws.Range("F10").Value = 2.774198266425014
$ws.Range("B11").Value = 0.02680902062944027
$ws.Range("C11").Value = 1.376487127339883
$ws.Range("D11").Value = 7.472110267386162
$ws.Range("E11").Value = 2.733516099712267
$ws.Range("F11").Value = 2.794816675622101
